$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D-column cells whose new values would otherwise be
# auto-parsed by Excel as numbers (these display as dotted "thousand.cents"
# style text, e.g. "249.77" must stay literal text "249.77", not a number).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "37.162.11"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.054.41"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "249.77"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").Value = "59.17"
$ws.Range("E7").Value = "  +6.73%  "
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "16.13"
$ws.Range("E12").Value = "  +6.31%  "
$ws.Range("D13").Value = "2.354.60"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "0.815"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "5.56"
$ws.Range("E15").Value = "  +5.84%  "
$ws.Range("D16").Value = "2.056.88"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "37.180.94"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "17.28"
$ws.Range("E18").Value = "  +21.82%  "
$ws.Range("D19").Value = "74.81"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "5.38"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "236.97"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  +11.16%  "
$ws.Range("D26").Value = "168.32"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "9.32"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").Value = "19.93"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +9.47%  "
$ws.Range("D31").Value = "4.74"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  +4.37%  "
$ws.Range("D34").Value = "0.0897"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  +6.32%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("E40").Value = "  +14.18%  "
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  +30.33%  "
$ws.Range("D42").Value = "17.62"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "0.0222"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "96.31"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "1.285.57"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "6.79"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "2.243.73"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -17.45%  "

# Restore default (unstyled) cell style for the cells we temporarily set to
# text format, so only the displayed values change -- matching source formatting.
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
